$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns (I, J)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header style used by the existing header cells (e.g. H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for columns I (I0) and J (IF), keyed by row number
$data = @"
2	8	8
3	8	8
4	5	6
5	6	7
6	5	5
7	7	7
8	6	6
9	6	6
10	8	9
11	7	8
12	7	7
13	6	7
14	4	5
15	7	7
16	4	5
17	6	6
18	6	7
19	5	6
20	8	9
21	4	6
22	9	9
23	5	7
24	7	8
25	6	7
26	8	9
27	8	9
28	8	9
29	8	8
30	6	6
31	8	8
32	9	9
33	7	7
34	7	7
35	5	6
36	8	8
37	7	8
38	7	7
39	6	7
40	8	8
41	8	9
42	7	7
43	5	6
44	8	8
45	6	7
46	8	8
47	7	8
48	8	8
49	9	9
50	6	7
51	5	6
52	6	7
53	8	9
54	9	9
55	9	9
56	7	8
57	7	8
58	5	6
59	7	8
60	7	8
61	7	8
62	8	8
63	5	7
64	6	8
65	8	9
66	8	8
67	8	8
68	8	9
69	6	7
70	4	5
71	7	7
72	9	9
73	9	9
74	8	8
75	9	9
76	8	8
77	4	4
"@

$lines = $data -split "`n"
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line -eq "") { continue }
    $parts = $line -split "`t"
    $row = [int]$parts[0]
    $iVal = [int]$parts[1]
    $jVal = [int]$parts[2]
    $ws.Cells.Item($row, 9).Value = $iVal
    $ws.Cells.Item($row, 10).Value = $jVal
}
